# Auto-generated edit script applying numeric updates to the Golem_Profits workbook
# (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 59.57143
$ws.Cells.Item(8, 9).Value = 56.833332
$ws.Cells.Item(8, 11).Value = 170.499996
$ws.Cells.Item(8, 13).Value = -31.49999600000001

$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).ClearContents()
$ws.Cells.Item(40, 14).Value = 0

$ws.Cells.Item(48, 8).Value = 1750
$ws.Cells.Item(48, 10).Value = 1500
$ws.Cells.Item(48, 12).Value = 4500
$ws.Cells.Item(48, 14).Value = -5084

$ws.Cells.Item(56, 8).Value = 1750
$ws.Cells.Item(56, 10).Value = 1500
$ws.Cells.Item(56, 12).Value = 4500
$ws.Cells.Item(56, 14).Value = -5568

$ws.Cells.Item(94, 8).Value = 1826.6666
$ws.Cells.Item(94, 9).Value = 1826.6666
$ws.Cells.Item(94, 11).Value = 1826.6666
$ws.Cells.Item(94, 13).Value = -1375.6666

$ws.Cells.Item(138, 8).Value = 3266
$ws.Cells.Item(138, 9).Value = 2500
$ws.Cells.Item(138, 10).Value = 3457.5
$ws.Cells.Item(138, 11).Value = 7500
$ws.Cells.Item(138, 12).Value = 10372.5
$ws.Cells.Item(138, 13).Value = -2360
$ws.Cells.Item(138, 14).Value = -20652.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2930.9333
$ws.Cells.Item(45, 9).Value = 1774.5555
$ws.Cells.Item(45, 11).Value = 1774.5555
$ws.Cells.Item(45, 13).Value = -1397.5555

$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).ClearContents()
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = 0

$ws.Cells.Item(96, 8).Value = 37782.43
$ws.Cells.Item(96, 10).Value = 37782.43
$ws.Cells.Item(96, 12).Value = 37782.43
$ws.Cells.Item(96, 14).Value = -43274.43

$ws.Cells.Item(132, 8).Value = 2076.2727
$ws.Cells.Item(132, 9).Value = 1783.9
$ws.Cells.Item(132, 11).Value = 5351.700000000001
$ws.Cells.Item(132, 13).Value = -2821.700000000001

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).ClearContents()
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1508
$ws.Cells.Item(134, 9).Value = 1550.3334
$ws.Cells.Item(134, 11).Value = 4651.0002
$ws.Cells.Item(134, 13).Value = -2116.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 35000
$ws.Cells.Item(41, 10).Value = 35000
$ws.Cells.Item(41, 12).Value = 35000
$ws.Cells.Item(41, 14).Value = -35856

$ws.Cells.Item(56, 8).Value = 40001
$ws.Cells.Item(56, 9).Value = 40001
$ws.Cells.Item(56, 11).Value = 40001
$ws.Cells.Item(56, 13).Value = -39156

$ws.Cells.Item(58, 8).Value = 524.4
$ws.Cells.Item(58, 9).Value = 524.4
$ws.Cells.Item(58, 11).Value = 524.4
$ws.Cells.Item(58, 13).Value = -321.4

$ws.Cells.Item(59, 8).Value = 54000
$ws.Cells.Item(59, 10).Value = 54000
$ws.Cells.Item(59, 12).Value = 54000
$ws.Cells.Item(59, 14).Value = -56290

$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).ClearContents()
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = 0

$ws.Cells.Item(99, 8).Value = 2500199.5
$ws.Cells.Item(99, 9).Value = 2500199.5
$ws.Cells.Item(99, 11).Value = 2500199.5
$ws.Cells.Item(99, 13).Value = -2498701.5

$ws.Cells.Item(124, 8).Value = 49985.715
$ws.Cells.Item(124, 10).Value = 49985.715
$ws.Cells.Item(124, 12).Value = 49985.715
$ws.Cells.Item(124, 14).Value = -54895.715

$ws.Cells.Item(126, 8).Value = 2500199.5
$ws.Cells.Item(126, 9).Value = 2500199.5
$ws.Cells.Item(126, 11).Value = 7500598.5
$ws.Cells.Item(126, 13).Value = -7498128.5

$ws.Cells.Item(132, 8).Value = 2161.625
$ws.Cells.Item(132, 9).Value = 2075.25
$ws.Cells.Item(132, 11).Value = 6225.75
$ws.Cells.Item(132, 13).Value = -3695.75

$ws.Cells.Item(136, 8).Value = 524.4
$ws.Cells.Item(136, 9).Value = 524.4
$ws.Cells.Item(136, 11).Value = 1573.2
$ws.Cells.Item(136, 13).Value = 976.8000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 712.25
$ws.Cells.Item(17, 9).Value = 599
$ws.Cells.Item(17, 11).Value = 1797
$ws.Cells.Item(17, 13).Value = -1628

$ws.Cells.Item(34, 8).Value = 13124.667
$ws.Cells.Item(34, 9).Value = 999
$ws.Cells.Item(34, 10).Value = 14227
$ws.Cells.Item(34, 11).Value = 2997
$ws.Cells.Item(34, 12).Value = 42681
$ws.Cells.Item(34, 13).Value = -2913
$ws.Cells.Item(34, 14).Value = -42849

$ws.Cells.Item(39, 8).Value = 23625
$ws.Cells.Item(39, 10).Value = 29833.334
$ws.Cells.Item(39, 12).Value = 89500.00199999999
$ws.Cells.Item(39, 14).Value = -90088.00199999999

$ws.Cells.Item(46, 8).Value = 733.3333
$ws.Cells.Item(46, 10).Value = 700
$ws.Cells.Item(46, 12).Value = 2100
$ws.Cells.Item(46, 14).Value = -2282

$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).ClearContents()
$ws.Cells.Item(49, 14).Value = 0

$ws.Cells.Item(51, 8).Value = 1495.75
$ws.Cells.Item(51, 9).Value = 1497.6666
$ws.Cells.Item(51, 10).Value = 1490
$ws.Cells.Item(51, 11).Value = 4492.9998
$ws.Cells.Item(51, 12).Value = 4470
$ws.Cells.Item(51, 13).Value = -4032.9998
$ws.Cells.Item(51, 14).Value = -5390

$ws.Cells.Item(55, 8).Value = 2042.0454
$ws.Cells.Item(55, 10).Value = 2251.6177
$ws.Cells.Item(55, 12).Value = 6754.853099999999
$ws.Cells.Item(55, 14).Value = -7108.853099999999

$ws.Cells.Item(64, 8).Value = 3110
$ws.Cells.Item(64, 9).Value = 2000
$ws.Cells.Item(64, 11).Value = 6000
$ws.Cells.Item(64, 13).Value = -5730

$ws.Cells.Item(67, 8).Value = 3110
$ws.Cells.Item(67, 9).Value = 2000
$ws.Cells.Item(67, 11).Value = 6000
$ws.Cells.Item(67, 13).Value = -5064

$ws.Cells.Item(132, 8).Value = 1952
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()

$ws.Cells.Item(139, 8).Value = 48799.816
$ws.Cells.Item(139, 9).Value = 4200
$ws.Cells.Item(139, 11).Value = 12600
$ws.Cells.Item(139, 13).Value = -7460

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).ClearContents()
$ws.Cells.Item(21, 14).Value = 0

$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 12).ClearContents()
$ws.Cells.Item(30, 14).Value = 0

$ws.Cells.Item(122, 8).Value = 2305
$ws.Cells.Item(122, 9).Value = 1027.2
$ws.Cells.Item(122, 11).Value = 3081.6
$ws.Cells.Item(122, 13).Value = -631.6000000000004

$ws.Cells.Item(123, 8).Value = 76499.75
$ws.Cells.Item(123, 10).Value = 76499.75
$ws.Cells.Item(123, 12).Value = 76499.75
$ws.Cells.Item(123, 14).Value = -81399.75

$ws.Cells.Item(132, 8).Value = 2520
$ws.Cells.Item(132, 9).Value = 1087.5
$ws.Cells.Item(132, 11).Value = 3262.5
$ws.Cells.Item(132, 13).Value = -732.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 7981.615
$ws.Cells.Item(16, 9).Value = 1251
$ws.Cells.Item(16, 10).Value = 45000
$ws.Cells.Item(16, 11).Value = 1251
$ws.Cells.Item(16, 12).Value = 45000
$ws.Cells.Item(16, 13).Value = -1081
$ws.Cells.Item(16, 14).Value = -45340

$ws.Cells.Item(32, 8).Value = 4822.875
$ws.Cells.Item(32, 9).Value = 1940.4286
$ws.Cells.Item(32, 11).Value = 1940.4286
$ws.Cells.Item(32, 13).Value = -1623.4286

$ws.Cells.Item(61, 8).Value = 967.6667
$ws.Cells.Item(61, 9).Value = 967.6667
$ws.Cells.Item(61, 11).Value = 967.6667
$ws.Cells.Item(61, 13).Value = -765.6667

$ws.Cells.Item(82, 8).Value = 1573.125
$ws.Cells.Item(82, 9).Value = 1375
$ws.Cells.Item(82, 10).Value = 1771.25
$ws.Cells.Item(82, 11).Value = 1375
$ws.Cells.Item(82, 12).Value = 1771.25
$ws.Cells.Item(82, 13).Value = -1014
$ws.Cells.Item(82, 14).Value = -2493.25

$ws.Cells.Item(85, 8).Value = 1573.125
$ws.Cells.Item(85, 9).Value = 1375
$ws.Cells.Item(85, 10).Value = 1771.25
$ws.Cells.Item(85, 11).Value = 1375
$ws.Cells.Item(85, 12).Value = 1771.25
$ws.Cells.Item(85, 13).Value = -127
$ws.Cells.Item(85, 14).Value = -4267.25

$ws.Cells.Item(113, 8).Value = 967.6667
$ws.Cells.Item(113, 9).Value = 967.6667
$ws.Cells.Item(113, 11).Value = 967.6667
$ws.Cells.Item(113, 13).Value = 1202.3333

$ws.Cells.Item(122, 8).Value = 1150
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 1595926.1
$ws.Cells.Item(136, 9).Value = 1004000.4
$ws.Cells.Item(136, 11).Value = 3012001.2
$ws.Cells.Item(136, 13).Value = -3009451.2

$ws.Cells.Item(137, 8).Value = 78000
$ws.Cells.Item(137, 10).Value = 120000
$ws.Cells.Item(137, 12).Value = 120000
$ws.Cells.Item(137, 14).Value = -130200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 49999.332
$ws.Cells.Item(124, 10).Value = 49999.5
$ws.Cells.Item(124, 12).Value = 49999.5
$ws.Cells.Item(124, 14).Value = -59819.5

$ws.Cells.Item(132, 8).Value = 899.3333
$ws.Cells.Item(132, 9).Value = 934.2857
$ws.Cells.Item(132, 11).Value = 2802.8571
$ws.Cells.Item(132, 13).Value = -272.8571000000002

$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).ClearContents()
$ws.Cells.Item(133, 14).Value = 0

$ws.Cells.Item(135, 8).Value = 71715
$ws.Cells.Item(135, 10).Value = 71715
$ws.Cells.Item(135, 12).Value = 71715
$ws.Cells.Item(135, 14).Value = -81855

$ws.Cells.Item(136, 8).Value = 1573.6666
$ws.Cells.Item(136, 9).Value = 1423.5264
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 4270.5792
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -1720.5792
$ws.Cells.Item(136, 14).Value = -14100

$ws.Cells.Item(137, 8).Value = 84500
$ws.Cells.Item(137, 10).Value = 84500
$ws.Cells.Item(137, 12).Value = 84500
$ws.Cells.Item(137, 14).Value = -94700

$ws.Cells.Item(139, 8).Value = 39999
$ws.Cells.Item(139, 10).Value = 39999
$ws.Cells.Item(139, 12).Value = 39999
$ws.Cells.Item(139, 14).Value = -50279

